$wb = $excel.ActiveWorkbook

# --- Update status to reflect a failed handback transform -----------------
$newStatus = "Handback transform failed"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = $newStatus
$overview.Range("C7").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = $newStatus

# --- Record the handback/handoff file-name mismatch in Error Detail -------
$zhcn.Range("L7").Value = "Handback file name: sfq3ovh4.kli is different with handoff file name: 7e22ed26-29e3-45c2-924d-9a71be0e7380.1256eab912a8893fc7285d96e25e8db050b3dea4.zh-cn."

$dede.Range("L7").Value = "Handback file name: sfq3ovh4.kli is different with handoff file name: 7e22ed26-29e3-45c2-924d-9a71be0e7380.1256eab912a8893fc7285d96e25e8db050b3dea4.de-de."
